# New weekly price record for Berenjena - Vega Monumental Concepción.
# A new row is inserted before row 111 (pushing existing rows 111-163 down
# to 112-164) and populated with the latest reported data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("111:111").Insert()

$ws.Range("A111").Value = 11
$ws.Range("B111").Value = "Vega Monumental Concepción"
$ws.Range("C111").Value = "Bíobío"
$ws.Range("D111").Value = 45097
$ws.Range("E111").Value = 8
$ws.Range("F111").Value = 100112001
$ws.Range("G111").Value = "Berenjena"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 220
$ws.Range("K111").Value = 6000
$ws.Range("L111").Value = 6500
$ws.Range("M111").Value = 6273
$ws.Range("N111").Value = "$/caja 50 unidades"
$ws.Range("O111").Value = "Región de Arica y Parinacota"
$ws.Range("P111").Value = 125
$ws.Range("Q111").Value = 50
$ws.Range("R111").Value = "Hortaliza"
